$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values recomputed for rows 2-11
$newValues = @{
    2  = 7
    3  = 3
    4  = 5
    5  = 4
    6  = 3
    7  = 3
    8  = 1
    9  = 6
    10 = 3
    11 = 3
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
